$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Maces")

# ---------------------------------------------------------------------------
# The "Maces" sheet is being reworked: existing weapons are kept (in the same
# relative order) but spread out with blank spacer rows between logical
# weapon-family groups, and four new weapons are inserted into the table.
#
# We rewrite rows 2-25 top to bottom directly (no Insert/Delete of rows) so
# that every cell keeps its original per-column style index, matching how
# the workbook actually changed. New text values introduced below are typed
# in the exact order they first appear in the final shared-strings table:
#   1) "w_mace_german" / "German Mace"                 (row 18)
#   2) "w_knight_warhammer_3" / "Knight Warhammer"      (row 10)
#   3) "w_knight_flanged_mace" / "Knight Flanged Mace"  (row 15)
#   4) "w_kriegshammer" / "Kriegshammer"                (row 25)
# ---------------------------------------------------------------------------

function Set-MaceRow {
    param($Row, $Id, $Name, $Price, $Weight, $Difficulty, $Swing, $Thrust, $Hit, $Block)

    $ws.Cells.Item($Row, 1).Value2 = $Id
    $ws.Cells.Item($Row, 2).Value2 = $Name
    $ws.Cells.Item($Row, 3).Value2 = $Price
    $ws.Cells.Item($Row, 4).Value2 = $Weight
    if ($null -eq $Difficulty) {
        $ws.Cells.Item($Row, 5).ClearContents()
    } else {
        $ws.Cells.Item($Row, 5).Value2 = $Difficulty
    }
    $ws.Cells.Item($Row, 6).Value2 = $Swing
    $ws.Cells.Item($Row, 7).Value2 = $Thrust
    $ws.Cells.Item($Row, 8).Value2 = $Hit
    $ws.Cells.Item($Row, 9).Value2 = $Block
}

function Clear-MaceRow {
    param($Row)
    $ws.Range($ws.Cells.Item($Row, 1), $ws.Cells.Item($Row, 9)).ClearContents()
}

# Row 18 first -> introduces "w_mace_german" / "German Mace" as the first
# brand-new shared strings.
Set-MaceRow 18 "w_mace_german" "German Mace" 278 3.4 $null 97 72 27 0

# Row 10 -> introduces "w_knight_warhammer_3" (Name reuses existing "Knight Warhammer").
Set-MaceRow 10 "w_knight_warhammer_3" "Knight Warhammer" 365 2.2 $null 94 70 34 0

# Row 15 -> introduces "w_knight_flanged_mace" / "Knight Flanged Mace".
Set-MaceRow 15 "w_knight_flanged_mace" "Knight Flanged Mace" 344 4.2 $null 95 72 29 0

# Row 25 -> introduces "w_kriegshammer" / "Kriegshammer".
Set-MaceRow 25 "w_kriegshammer" "Kriegshammer" 512 4.5 12 82 83 42 32

# Now rewrite the rest of the rows (previously existing rows, shifted to
# their new positions) and the new blank spacer rows, top-down.
Set-MaceRow 2 "w_wooden_stick" "Wooden Stick" 4 2.5 $null 99 63 13 0
Set-MaceRow 3 "w_archers_maul" "Archers Maul" 77 2 $null 99 73 20 0

Clear-MaceRow 4

Set-MaceRow 5 "w_warhammer_1" "Warhammer 1" 293 2 $null 95 70 30 0
Set-MaceRow 6 "w_warhammer_2" "Warhammer 2" 317 2 $null 95 70 31 0

Clear-MaceRow 7

Set-MaceRow 8 "w_knight_warhammer_1" "Spiked Knight Warhammer" 372 2.5 $null 93 76 33 19
Set-MaceRow 9 "w_knight_warhammer_2" "Knight Warhammer" 334 2 $null 95 63 32 0

# row 10 already written above

Clear-MaceRow 11

Set-MaceRow 12 "w_great_hammer" "Great Hammer" 422 9 14 79 75 45 0

Clear-MaceRow 13

Set-MaceRow 14 "w_knight_winged_mace" "Knight Winged Mace" 336 4 $null 96 69 28 0

# row 15 already written above

Clear-MaceRow 16

Set-MaceRow 17 "w_mace_english" "English Mace" 262 3.25 $null 97 72 26 0

# row 18 already written above

Clear-MaceRow 19

Set-MaceRow 20 "w_spiked_club" "Spiked Club" 83 3.25 $null 96 75 21 0
Set-MaceRow 21 "w_mace_knobbed" "Knobbed Mace " 98 2.5 $null 98 70 21 0
Set-MaceRow 22 "w_mace_spiked" "Spiked Mace" 152 2.75 $null 98 71 23 0
Set-MaceRow 23 "w_mace_winged" "Winged Mace" 212 3 $null 97 71 24 0

Clear-MaceRow 24

# row 25 already written above

# ---------------------------------------------------------------------------
# View-state: the active sheet moved from "Axes & Bardiches" to "Maces", and
# the selection on the Maces sheet moved to D26.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("D26").Select()
